$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A120:D120").Copy($ws.Range("A121:D121"))
$ws.Range("A121").Value2 = 45140
$ws.Range("B121").Value2 = 0.54917824074074073
$ws.Range("C121").Value2 = 83247
$ws.Range("D121").Value2 = 1690

$ws.Range("A121:D121").Copy($ws.Range("A122:D122"))
$ws.Range("A122").Value2 = 45141
$ws.Range("B122").Value2 = 0.64541666666666664
$ws.Range("C122").Value2 = 83254
$ws.Range("D122").Value2 = 1690

$ws.Range("A122:D122").Copy($ws.Range("A123:D123"))
$ws.Range("A123").Value2 = 45142
$ws.Range("B123").Value2 = 0.43209490740740741
$ws.Range("C123").Value2 = 83255
$ws.Range("D123").Value2 = 1690

$ws.Range("A123:D123").Copy($ws.Range("A124:D124"))
$ws.Range("A124").Value2 = 45143
$ws.Range("B124").Value2 = 0.47541666666666665
$ws.Range("C124").Value2 = 83255
$ws.Range("D124").Value2 = 1690

$ws.Range("A124:D124").Copy($ws.Range("A125:D125"))
$ws.Range("A125").Value2 = 45144
$ws.Range("B125").Value2 = 0.4729976851851852
$ws.Range("C125").Value2 = 83255
$ws.Range("D125").Value2 = 1690

$ws.Range("A125:D125").Copy($ws.Range("A126:D126"))
$ws.Range("A126").Value2 = 45145
$ws.Range("B126").Value2 = 0.4765625
$ws.Range("C126").Value2 = 83280
$ws.Range("D126").Value2 = 1690

$ws.Range("A126:D126").Copy($ws.Range("A127:D127"))
$ws.Range("A127").Value2 = 45146
$ws.Range("B127").Value2 = 0.46464120370370371
$ws.Range("C127").Value2 = 83280
$ws.Range("D127").Value2 = 1690

$ws.Range("A127:D127").Copy($ws.Range("A128:D128"))
$ws.Range("A128").Value2 = 45147
$ws.Range("B128").Value2 = 0.36782407407407408
$ws.Range("C128").Value2 = 83284
$ws.Range("D128").Value2 = 1690

$ws.Range("A128:D128").Copy($ws.Range("A129:D129"))
$ws.Range("A129").Value2 = 45148
$ws.Range("B129").Value2 = 0.48040509259259262
$ws.Range("C129").Value2 = 83287
$ws.Range("D129").Value2 = 1690

$ws.Range("A129:D129").Copy($ws.Range("A130:D130"))
$ws.Range("A130").Value2 = 45149
$ws.Range("B130").Value2 = 0.52964120370370371
$ws.Range("C130").Value2 = 83291
$ws.Range("D130").Value2 = 1690

$ws.Range("A130:D130").Copy($ws.Range("A131:D131"))
$ws.Range("A131").Value2 = 45150
$ws.Range("B131").Value2 = 0.47501157407407407
$ws.Range("C131").Value2 = 83359
$ws.Range("D131").Value2 = 1690

$ws.Range("A131:D131").Copy($ws.Range("A132:D132"))
$ws.Range("A132").Value2 = 45151
$ws.Range("B132").Value2 = 0.46175925925925926
$ws.Range("C132").Value2 = 83359
$ws.Range("D132").Value2 = 1690

$ws.Range("A132:D132").Copy($ws.Range("A133:D133"))
$ws.Range("A133").Value2 = 45152
$ws.Range("B133").Value2 = 0.46082175925925922
$ws.Range("C133").Value2 = 83359
$ws.Range("D133").Value2 = 1690

$ws.Range("A133:D133").Copy($ws.Range("A134:D134"))
$ws.Range("A134").Value2 = 45153
$ws.Range("B134").Value2 = 0.39368055555555559
$ws.Range("C134").Value2 = 83359
$ws.Range("D134").Value2 = 1690

$ws.Range("A134:D134").Copy($ws.Range("A135:D135"))
$ws.Range("A135").Value2 = 45154
$ws.Range("B135").Value2 = 0.41212962962962968
$ws.Range("C135").Value2 = 83360
$ws.Range("D135").Value2 = 1690

$ws.Range("A135:D135").Copy($ws.Range("A136:D136"))
$ws.Range("A136").Value2 = 45155
$ws.Range("B136").Value2 = 0.42332175925925924
$ws.Range("C136").Value2 = 77977
$ws.Range("D136").Value2 = 1680

$ws.Range("A137").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 116